$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update progress values in column F
$ws.Range("F5").Value = 0.8
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 0.25
$ws.Range("F11").Value = 0.25

# Remove the "Demo" task row (row 15) entirely, which shifts nothing below it
$ws.Rows.Item(15).Delete()
